# Auto-generated edit script for 세아홀딩스 IFRS workbook
# Commit message: "error solve ifrs list"
# Updates the financial figures for annual rows (2014-2018, rows 2-6)
# and clears the erroneous forecast estimate rows (2019E-2021E, rows 7-9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 38961
$ws.Range("E2").Value = 2556
$ws.Range("F2").Value = 2655
$ws.Range("G2").Value = 2355
$ws.Range("H2").Value = 1645
$ws.Range("I2").Value = 974
$ws.Range("J2").Value = 671
$ws.Range("K2").Value = 41411
$ws.Range("L2").Value = 18243
$ws.Range("M2").Value = 23167
$ws.Range("N2").Value = 14865
$ws.Range("O2").Value = 8302
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 2946
$ws.Range("R2").Value = -1058
$ws.Range("S2").Value = 652
$ws.Range("T2").Value = 945
$ws.Range("U2").Value = 2001
$ws.Range("V2").Value = 13112
$ws.Range("W2").Value = 6.56
$ws.Range("X2").Value = 4.22
$ws.Range("Y2").Value = 6.76
$ws.Range("Z2").Value = 4.08
$ws.Range("AA2").Value = 78.75
$ws.Range("AB2").Value = 7334.96
$ws.Range("AC2").Value = 24345
$ws.Range("AD2").Value = 6.65
$ws.Range("AE2").Value = 371715
$ws.Range("AF2").Value = 0.44
$ws.Range("AG2").Value = 1800
$ws.Range("AH2").Value = 1.11
$ws.Range("AI2").Value = 7.39
$ws.Range("AJ2").Value = 4000000

# --- Row 3 ---
$ws.Range("D3").Value = 40482
$ws.Range("E3").Value = 2876
$ws.Range("F3").Value = 2876
$ws.Range("G3").Value = 2649
$ws.Range("H3").Value = 2071
$ws.Range("I3").Value = 1159
$ws.Range("J3").Value = 911
$ws.Range("K3").Value = 51915
$ws.Range("L3").Value = 24038
$ws.Range("M3").Value = 27877
$ws.Range("N3").Value = 15260
$ws.Range("O3").Value = 12617
$ws.Range("P3").Value = 200
$ws.Range("Q3").Value = 5009
$ws.Range("R3").Value = -4807
$ws.Range("S3").Value = -1639
$ws.Range("T3").Value = 1194
$ws.Range("U3").Value = 3815
$ws.Range("V3").Value = 14861
$ws.Range("W3").Value = 7.1
$ws.Range("X3").Value = 5.12
$ws.Range("Y3").Value = 7.7
$ws.Range("Z3").Value = 4.44
$ws.Range("AA3").Value = 86.23
$ws.Range("AB3").Value = 7469.51
$ws.Range("AC3").Value = 28983
$ws.Range("AD3").Value = 4.52
$ws.Range("AE3").Value = 381583
$ws.Range("AF3").Value = 0.34
$ws.Range("AG3").Value = 1750
$ws.Range("AH3").Value = 1.34
$ws.Range("AI3").Value = 6.04
$ws.Range("AJ3").Value = 4000000

# --- Row 4 ---
$ws.Range("D4").Value = 40343
$ws.Range("E4").Value = 2348
$ws.Range("F4").Value = 2348
$ws.Range("G4").Value = 2031
$ws.Range("H4").Value = 1619
$ws.Range("I4").Value = 901
$ws.Range("J4").Value = 718
$ws.Range("K4").Value = 51779
$ws.Range("L4").Value = 22537
$ws.Range("M4").Value = 29241
$ws.Range("N4").Value = 17185
$ws.Range("O4").Value = 12057
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 3825
$ws.Range("R4").Value = -2236
$ws.Range("S4").Value = -2216
$ws.Range("T4").Value = 1974
$ws.Range("U4").Value = 1851
$ws.Range("V4").Value = 14623
$ws.Range("W4").Value = 5.82
$ws.Range("X4").Value = 4.01
$ws.Range("Y4").Value = 5.55
$ws.Range("Z4").Value = 3.12
$ws.Range("AA4").Value = 77.06999999999999
$ws.Range("AB4").Value = 8458.34
$ws.Range("AC4").Value = 22514
$ws.Range("AD4").Value = 5.73
$ws.Range("AE4").Value = 429719
$ws.Range("AF4").Value = 0.3
$ws.Range("AG4").Value = 1750
$ws.Range("AH4").Value = 1.36
$ws.Range("AI4").Value = 7.77
$ws.Range("AJ4").Value = 4000000

# --- Row 5 ---
$ws.Range("D5").Value = 47944
$ws.Range("E5").Value = 2746
$ws.Range("F5").Value = 2746
$ws.Range("G5").Value = 2492
$ws.Range("H5").Value = 2111
$ws.Range("I5").Value = 1344
$ws.Range("J5").Value = 767
$ws.Range("K5").Value = 52548
$ws.Range("L5").Value = 22088
$ws.Range("M5").Value = 30460
$ws.Range("N5").Value = 18678
$ws.Range("O5").Value = 11782
$ws.Range("P5").Value = 200
$ws.Range("Q5").Value = 2647
$ws.Range("R5").Value = -1641
$ws.Range("S5").Value = -1473
$ws.Range("T5").Value = 1288
$ws.Range("U5").Value = 1359
$ws.Range("V5").Value = 13838
$ws.Range("W5").Value = 5.73
$ws.Range("X5").Value = 4.4
$ws.Range("Y5").Value = 7.5
$ws.Range("Z5").Value = 4.05
$ws.Range("AA5").Value = 72.52
$ws.Range("AB5").Value = 9253.73
$ws.Range("AC5").Value = 33604
$ws.Range("AD5").Value = 4.52
$ws.Range("AE5").Value = 467055
$ws.Range("AF5").Value = 0.33
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1.32
$ws.Range("AI5").Value = 5.95
$ws.Range("AJ5").Value = 4000000

# --- Row 6 ---
$ws.Range("D6").Value = 51769
$ws.Range("E6").Value = 1175
$ws.Range("F6").Value = 1175
$ws.Range("G6").Value = 783
$ws.Range("H6").Value = 645
$ws.Range("I6").Value = 328
$ws.Range("K6").Value = 54539
$ws.Range("L6").Value = 24466
$ws.Range("M6").Value = 30073
$ws.Range("N6").Value = 19386
$ws.Range("P6").Value = 200
$ws.Range("Q6").Value = 1407
$ws.Range("R6").Value = -1302
$ws.Range("S6").Value = 655
$ws.Range("T6").Value = 1177
$ws.Range("U6").Value = 230
$ws.Range("V6").Value = 16106
$ws.Range("W6").Value = 2.27
$ws.Range("X6").Value = 1.25
$ws.Range("Y6").Value = 1.72
$ws.Range("Z6").Value = 1.21
$ws.Range("AA6").Value = 81.36
$ws.Range("AB6").Value = 9645.049999999999
$ws.Range("AC6").Value = 8191
$ws.Range("AD6").Value = 11.88
$ws.Range("AE6").Value = 484760
$ws.Range("AF6").Value = 0.2
$ws.Range("AG6").Value = 2500
$ws.Range("AH6").Value = 2.57
$ws.Range("AI6").Value = 30.51
$ws.Range("AJ6").Value = 4000000

# Forecast rows (2019E/2020E/2021E) had erroneous data; clear all value cells,
# keeping only the row index (A), period label (B) and year label (C).
$ws.Range("D7:AI9").ClearContents()

Write-Host "Applied ifrs list fix for rows 2-9"
